$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3640.9167
$ws.Range("I64").Value = 3426.5454
$ws.Range("K64").Value = 3426.5454
$ws.Range("M64").Value = -3178.5454
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 3640.9167
$ws.Range("I67").Value = 3426.5454
$ws.Range("K67").Value = 3426.5454
$ws.Range("M67").Value = -2568.5454
$ws.Range("N67").ClearContents()

$ws.Range("H70").Value = 1451.5238
$ws.Range("I70").Value = 1610.2
$ws.Range("J70").Value = 1307.2727
$ws.Range("K70").Value = 4830.6
$ws.Range("L70").Value = 3921.8181
$ws.Range("M70").Value = -4560.6
$ws.Range("N70").Value = -4461.8181

$ws.Range("H73").Value = 1451.5238
$ws.Range("I73").Value = 1610.2
$ws.Range("J73").Value = 1307.2727
$ws.Range("K73").Value = 4830.6
$ws.Range("L73").Value = 3921.8181
$ws.Range("M73").Value = -3894.6
$ws.Range("N73").Value = -5793.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1337806.2
$ws.Range("I32").Value = 1523410.8
$ws.Range("J32").Value = 106067.45
$ws.Range("K32").Value = 1523410.8
$ws.Range("L32").Value = 106067.45
$ws.Range("M32").Value = -1523123.8
$ws.Range("N32").Value = -106641.45

$ws.Range("H53").Value = 65042
$ws.Range("I53").Value = 20039
$ws.Range("J53").Value = 80043
$ws.Range("K53").Value = 20039
$ws.Range("L53").Value = 80043
$ws.Range("M53").Value = -19357
$ws.Range("N53").Value = -81407

$ws.Range("H102").Value = 2711.2856
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 2829.8333
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 2829.8333
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -6073.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 925.6875
$ws.Range("I20").Value = 860.5
$ws.Range("J20").Value = 990.875
$ws.Range("K20").Value = 860.5
$ws.Range("L20").Value = 990.875
$ws.Range("M20").Value = -613.5
$ws.Range("N20").Value = -1484.875

$ws.Range("H105").Value = 2000
$ws.Range("I105").Value = 1850
$ws.Range("J105").Value = 2150
$ws.Range("K105").Value = 1850
$ws.Range("L105").Value = 2150
$ws.Range("M105").Value = -103
$ws.Range("N105").Value = -5644

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 83335600
$ws.Range("I62").Value = 83335600
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 83335600
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -83334976
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 83335600
$ws.Range("I65").Value = 83335600
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 416678000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -416674880
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 6282
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 6282
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 18846
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -20218

$ws.Range("H63").Value = 3014
$ws.Range("I63").Value = 3014
$ws.Range("J63").Value = 3014
$ws.Range("K63").Value = 9042
$ws.Range("L63").Value = 9042
$ws.Range("M63").Value = -8293
$ws.Range("N63").Value = -10540

$ws.Range("H64").Value = 3793104.2
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3793104.2
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 11379312.6
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -11379852.6

$ws.Range("H65").Value = 6282
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 6282
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 56538
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -63402

$ws.Range("H66").Value = 3014
$ws.Range("I66").Value = 3014
$ws.Range("J66").Value = 3014
$ws.Range("K66").Value = 27126
$ws.Range("L66").Value = 27126
$ws.Range("M66").Value = -23382
$ws.Range("N66").Value = -34614

$ws.Range("H67").Value = 3793104.2
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3793104.2
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 11379312.6
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -11381184.6

$ws.Range("H68").Value = 2271.1265
$ws.Range("I68").Value = 4720.625
$ws.Range("J68").Value = 1202.2545
$ws.Range("K68").Value = 14161.875
$ws.Range("L68").Value = 3606.7635
$ws.Range("M68").Value = -13350.875
$ws.Range("N68").Value = -5228.7635

$ws.Range("H69").Value = 33335424
$ws.Range("J69").Value = 33335424
$ws.Range("L69").Value = 100006272
$ws.Range("N69").Value = -100007894

$ws.Range("H70").Value = 2787.3
$ws.Range("I70").Value = 1029.6
$ws.Range("J70").Value = 3373.2
$ws.Range("K70").Value = 3088.8
$ws.Range("L70").Value = 10119.6
$ws.Range("M70").Value = -2773.8
$ws.Range("N70").Value = -10749.6

$ws.Range("H71").Value = 2271.1265
$ws.Range("I71").Value = 4720.625
$ws.Range("J71").Value = 1202.2545
$ws.Range("K71").Value = 42485.625
$ws.Range("L71").Value = 10820.2905
$ws.Range("M71").Value = -38429.625
$ws.Range("N71").Value = -18932.2905

$ws.Range("H72").Value = 33335424
$ws.Range("J72").Value = 33335424
$ws.Range("L72").Value = 300018816
$ws.Range("N72").Value = -300026928

$ws.Range("H73").Value = 2787.3
$ws.Range("I73").Value = 1029.6
$ws.Range("J73").Value = 3373.2
$ws.Range("K73").Value = 3088.8
$ws.Range("L73").Value = 10119.6
$ws.Range("M73").Value = -1996.8
$ws.Range("N73").Value = -12303.6

$ws.Range("H74").Value = 3048.5386
$ws.Range("I74").Value = 1513
$ws.Range("J74").Value = 3176.5
$ws.Range("K74").Value = 4539
$ws.Range("L74").Value = 9529.5
$ws.Range("M74").Value = -3478
$ws.Range("N74").Value = -11651.5

$ws.Range("H75").Value = 6738.3335
$ws.Range("J75").Value = 6763.421
$ws.Range("L75").Value = 20290.263
$ws.Range("N75").Value = -22286.263

$ws.Range("H76").Value = 10074.846
$ws.Range("I76").Value = 3670.6667
$ws.Range("J76").Value = 10910.174
$ws.Range("K76").Value = 11012.0001
$ws.Range("L76").Value = 32730.522
$ws.Range("M76").Value = -10629.0001
$ws.Range("N76").Value = -33496.522

$ws.Range("H77").Value = 3048.5386
$ws.Range("I77").Value = 1513
$ws.Range("J77").Value = 3176.5
$ws.Range("K77").Value = 13617
$ws.Range("L77").Value = 28588.5
$ws.Range("M77").Value = -8313
$ws.Range("N77").Value = -39196.5

$ws.Range("H78").Value = 6738.3335
$ws.Range("J78").Value = 6763.421
$ws.Range("L78").Value = 60870.789
$ws.Range("N78").Value = -70854.789

$ws.Range("H79").Value = 10074.846
$ws.Range("I79").Value = 3670.6667
$ws.Range("J79").Value = 10910.174
$ws.Range("K79").Value = 11012.0001
$ws.Range("L79").Value = 32730.522
$ws.Range("M79").Value = -9686.000100000001
$ws.Range("N79").Value = -35382.522

$ws.Range("H82").Value = 200003200
$ws.Range("I82").Value = 1013
$ws.Range("J82").Value = 250003740
$ws.Range("K82").Value = 3039
$ws.Range("L82").Value = 750011220
$ws.Range("M82").Value = -2633
$ws.Range("N82").Value = -750012032

$ws.Range("H85").Value = 200003200
$ws.Range("I85").Value = 1013
$ws.Range("J85").Value = 250003740
$ws.Range("K85").Value = 3039
$ws.Range("L85").Value = 750011220
$ws.Range("M85").Value = -1635
$ws.Range("N85").Value = -750014028

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4250
$ws.Range("I70").Value = 4080
$ws.Range("J70").Value = 4344.4443
$ws.Range("K70").Value = 4080
$ws.Range("L70").Value = 4344.4443
$ws.Range("M70").Value = -3810
$ws.Range("N70").Value = -4884.4443

$ws.Range("H73").Value = 4250
$ws.Range("I73").Value = 4080
$ws.Range("J73").Value = 4344.4443
$ws.Range("K73").Value = 4080
$ws.Range("L73").Value = 4344.4443
$ws.Range("M73").Value = -3144
$ws.Range("N73").Value = -6216.4443

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()

$ws.Range("H93").Value = 20500
$ws.Range("J93").Value = 20500
$ws.Range("L93").Value = 20500
$ws.Range("N93").Value = -25492

$ws.Range("H138").Value = 54429
$ws.Range("J138").Value = 54429
$ws.Range("L138").Value = 54429
$ws.Range("N138").Value = -64709
